$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number + report week dates) ---
$ws.Range("A8").Characters(21, 2).Text = "51"
$ws.Range("C9").Characters(27, 9).Text = "12/16/2024"
$ws.Range("C9").Characters(48, 10).Text = "12/22/2024"

# --- Crime-data cell updates (rows 16-30) ---
$ws.Range("C14").Copy()
$ws.Range("C16").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("D16").Value = 2
$ws.Range("J14").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = -28.571428571428
$ws.Range("J16").Value = 66
$ws.Range("K16").Value = -9.090909090909
$ws.Range("L16").Value = -6.25
$ws.Range("M16").Value = -40
$ws.Range("N16").Value = -84.575835475578
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = -11.111111111111
$ws.Range("I17").Value = 119
$ws.Range("J17").Value = 102
$ws.Range("K17").Value = 16.666666666666
$ws.Range("L17").Value = 16.666666666666
$ws.Range("M17").Value = 32.222222222222
$ws.Range("N17").Value = -60.333333333333
$ws.Range("C18").Value = 2
$ws.Range("C14").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E18").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = 200
$ws.Range("I18").Value = 85
$ws.Range("K18").Value = 3.658536585365
$ws.Range("L18").Value = -23.423423423423
$ws.Range("M18").Value = -18.269230769230
$ws.Range("N18").Value = -80.046948356807
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 300
$ws.Range("F19").Value = 19
$ws.Range("G19").Value = 11
$ws.Range("H19").Value = 72.727272727272
$ws.Range("I19").Value = 205
$ws.Range("J19").Value = 162
$ws.Range("K19").Value = 26.543209876543
$ws.Range("L19").Value = -3.301886792452
$ws.Range("M19").Value = -9.691629955947
$ws.Range("N19").Value = -10.869565217391
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = -80
$ws.Range("L20").Value = 7.407407407407
$ws.Range("M20").Value = -1.694915254237
$ws.Range("N20").Value = -85.316455696202
$ws.Range("C21").Value = 8
$ws.Range("D21").Value = 7
$ws.Range("E21").Value = 14.285714285714
$ws.Range("F21").Value = 39
$ws.Range("G21").Value = 34
$ws.Range("H21").Value = 14.705882352941
$ws.Range("I21").Value = 532
$ws.Range("J21").Value = 480
$ws.Range("K21").Value = 10.833333333333
$ws.Range("L21").Value = -3.623188405797
$ws.Range("M21").Value = -8.747855917667
$ws.Range("N21").Value = -69.772727272727
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 6
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 15
$ws.Range("H23").Value = -53.333333333333
$ws.Range("I23").Value = 106
$ws.Range("J23").Value = 111
$ws.Range("K23").Value = -4.504504504504
$ws.Range("L23").Value = -5.357142857142
$ws.Range("M23").Value = 34.177215189873
$ws.Range("C24").Value = 7
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = -22.222222222222
$ws.Range("F24").Value = 38
$ws.Range("G24").Value = 28
$ws.Range("H24").Value = 35.714285714285
$ws.Range("I24").Value = 525
$ws.Range("J24").Value = 579
$ws.Range("K24").Value = -9.326424870466
$ws.Range("L24").Value = 17.977528089887
$ws.Range("M24").Value = 22.950819672131
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 0
$ws.Range("G25").Value = 10
$ws.Range("H25").Value = 30
$ws.Range("I25").Value = 252
$ws.Range("J25").Value = 294
$ws.Range("K25").Value = -14.285714285714
$ws.Range("L25").Value = 100
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 200
$ws.Range("F26").Value = 9
$ws.Range("G26").Value = 13
$ws.Range("H26").Value = -30.769230769230
$ws.Range("I26").Value = 162
$ws.Range("J26").Value = 152
$ws.Range("K26").Value = 6.578947368421
$ws.Range("L26").Value = -2.409638554216
$ws.Range("M26").Value = -42.553191489361
$ws.Range("C14").Copy()
$ws.Range("G28").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("G28").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("H28").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("H28").PasteSpecial(-4122)
$ws.Range("G29").Value = 3
$ws.Range("J29").Value = 11
$ws.Range("K29").Value = -63.636363636363
$ws.Range("G30").Value = 3
$ws.Range("J30").Value = 11
$ws.Range("K30").Value = -63.636363636363

# --- Row insert: shift old rows 56-57 down to 57-58 ---
$ws.Rows.Item(56).Insert()
$ws.Range("A56").Clear()
